# Drinkable.pptx edit:
#  1. Insert a new "Title and Content" slide at position 3 with
#     Title = "MVP" and body = "Just to display a list of cocktails
#     based on ingredient choice."
#  2. Move the "Planning" slide (originally slide 5, position 6 once
#     the new slide has been inserted) up to position 4, right after
#     the new MVP slide. This pushes "Trials and Tribulations" and
#     "More Issues" back by one position each; everything else keeps
#     its relative order.

$p = $ppt.ActivePresentation

# --- 1. Insert the new MVP slide at position 3 ------------------------
# Layout 2 == ppLayoutText ("Title and Content"), matching the layout
# used by all the other text-content slides in this deck.
$mvpSlide = $p.Slides.Add(3, 2)
$mvpSlide.Shapes.Item(1).TextFrame.TextRange.Text = "MVP"
$mvpSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Just to display a list of cocktails based on ingredient choice."

# --- 2. Move "Planning" slide up to position 4 -------------------------
# Before this script runs the deck order is:
#   1 Drinkable
#   2 What does it do?
#   3 Trials and Tribulations
#   4 More Issues
#   5 Planning
#   6 Using inner-joins in SQLite
#   7 Using a spinner
#   8 Features to add
#   9 Things to do better.
#  10 Questions?
# After inserting the MVP slide at 3, "Planning" (originally 5) is now
# at position 6.
$planningSlide = $p.Slides.Item(6)
$planningSlide.MoveTo(4)
